$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "view the balance of a specific account" (row 4) as DONE
$ws.Range("C4").Value = "DONE"

# Update the selected cell to reflect the new active cell after the edit
$ws.Range("C4").Select()
